# Reorder the "Periodo Mora" values in column E (rows 16-19) from
# descending (2506,2505,2504,2503) to ascending (2503,2504,2505,2506)
# as part of updating the EC database with the new period batch.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2503"
$ws.Range("E17").Value = "2504"
$ws.Range("E18").Value = "2505"
$ws.Range("E19").Value = "2506"
